# Updated capital structure database
# Applies refreshed metrics to the Spain Brokerage & Investment Banking
# data rows (rows 2-4) in the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value = 0.21565
$ws.Range("E2").Value = -0.08905000000000002
$ws.Range("G2").Value = 0.0002865041467705454
$ws.Range("K2").Value = 48
$ws.Range("L2").Value = 0.1206333249560191
$ws.Range("M2").Value = 19.8528
$ws.Range("N2").Value = 0.0203264052421419
$ws.Range("O2").Value = 0.4136
$ws.Range("P2").Value = 19.8528
$ws.Range("Q2").Value = 0.0203264052421419
$ws.Range("R2").Value = 0.4136
$ws.Range("U2").Value = 1056.4
$ws.Range("V2").Value = 1.081601310535477
$ws.Range("W2").Value = 0.1431642494453109
$ws.Range("X2").Value = 0.03598227611209486
$ws.Range("Y2").Value = 0.107181973333216
$ws.Range("Z2").Value = -2.453144266337855
$ws.Range("AB2").Value = 0.03568152861580085
$ws.Range("AC2").Value = -0.03568152861580085
$ws.Range("AD2").Value = 32.14
$ws.Range("AF2").Value = 32.14
$ws.Range("AG2").Value = -1024.26
$ws.Range("AH2").Value = 0.0318583719915943
$ws.Range("AI2").Value = 0.06318810946838629
$ws.Range("AJ2").Value = 21.53616484440709
$ws.Range("AK2").Value = 1.869906528406601
$ws.Range("AM2").Value = -1.83

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value = 0.0503
$ws.Range("E3").Value = 0.0379
$ws.Range("K3").Value = 19.3
$ws.Range("L3").Value = 0.09984480082772891
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 961
$ws.Range("V3").Value = 2.733219567690557
$ws.Range("W3").Value = 0.1673894189071986
$ws.Range("X3").Value = 0.035628172419644
$ws.Range("Y3").Value = 0.1317612464875546
$ws.Range("Z3").Value = -0.640490390987409
$ws.Range("AB3").Value = 0.03552503595164139
$ws.Range("AC3").Value = -0.03552503595164139
$ws.Range("AD3").Value = 3.44
$ws.Range("AF3").Value = 3.44
$ws.Range("AG3").Value = -957.5599999999999
$ws.Range("AH3").Value = 0.009689049121225776
$ws.Range("AI3").Value = 0.02521254763998827
$ws.Range("AJ3").Value = 1.580236319228992
$ws.Range("AK3").Value = 1.161298146890463

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = 0.381
$ws.Range("E4").Value = -0.216
$ws.Range("G4").Value = 0.0005571847507331378
$ws.Range("K4").Value = 28.7
$ws.Range("L4").Value = 0.1402737047898338
$ws.Range("M4").Value = 19.8528
$ws.Range("N4").Value = 0.0317593984962406
$ws.Range("O4").Value = 0.6917351916376306
$ws.Range("P4").Value = 19.8528
$ws.Range("Q4").Value = 0.0317593984962406
$ws.Range("R4").Value = 0.6917351916376306
$ws.Range("U4").Value = 95.40000000000001
$ws.Range("V4").Value = 0.1526155815069589
$ws.Range("W4").Value = 0.1189390799834231
$ws.Range("X4").Value = 0.03633637980454572
$ws.Range("Y4").Value = 0.08260270017887739
$ws.Range("Z4").Value = 1.465616045845272
$ws.Range("AB4").Value = 0.03583802127996031
$ws.Range("AC4").Value = -0.03583802127996031
$ws.Range("AD4").Value = 28.7
$ws.Range("AF4").Value = 28.7
$ws.Range("AG4").Value = -66.7
$ws.Range("AH4").Value = 0.04389721627408993
$ws.Range("AI4").Value = 0.07710908113917249
$ws.Range("AJ4").Value = -0.1194484240687679
$ws.Range("AK4").Value = -0.2409682080924855
$ws.Range("AM4").Value = -1.83
